# Update the "cryptos" worksheet with refreshed price/volume figures
# (and a few re-ranked rows) as produced by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '51.720.53'
$ws.Range('D3').Value = '3.031.25'
$ws.Range('E3').Value = '  +2.49%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '380.78'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.10'
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.89'
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').Value = '3.509.16'
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.76'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = '3.028.20'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('E17').Value = '  -3.90%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.58'
$ws.Range('E18').Value = '  -16.92%  '
$ws.Range('D19').Value = '51.722.01'
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.54'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('E22').Value = '  +1.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.16'
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '268.69'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('E25').Value = '  -3.98%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.26'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.68'
$ws.Range('E27').Value = '  +9.42%  '
$ws.Range('E28').Value = '  +5.57%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.31'
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.11'
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '34.15'
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.52'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0449'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.31'
$ws.Range('E38').Value = '  +5.45%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.298'
$ws.Range('E39').Value = '  +15.31%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.05'
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.87'
$ws.Range('E41').Value = '  +2.48%  '
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '127.53'
$ws.Range('E43').Value = '  +6.24%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.116'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.78'
$ws.Range('E45').Value = '  +5.91%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '21.86'
$ws.Range('E46').Value = '  +1.97%  '
$ws.Range('E47').Value = '  +3.58%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  +4.90%  '
$ws.Range('D49').Value = '2.035.88'
$ws.Range('E49').Value = '  +1.29%  '
$ws.Range('D50').Value = '3.332.40'
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0321'
$ws.Range('E51').Value = '  +0.08%  '
